$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.390.95"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.572.10"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.78"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3762"
$ws.Range("E7").Value = "  +2.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.03"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3424"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07652"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.153"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.21"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.004"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.952"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "1.571.72"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.96"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06746"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.220"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "22.382.80"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.395"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.684"
$ws.Range("E26").Value = "  -9.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.22"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.10"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.030"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.32"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "1.746.15"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.165"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.010"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9837"
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.959"
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08516"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02552"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.393"
$ws.Range("E38").Value = "  +11.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2319"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06566"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.435"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6406"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.04"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.787"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5983"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.296"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.40"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07328"
$ws.Range("E51").Value = "  +0.64%  "
